$d = $word.ActiveDocument

# Helper: replace the first occurrence of $find inside paragraph $p's range
# with $replace, and force the replaced span to become its own run(s) by
# toggling a character-formatting property on and back off again (this is
# how Word naturally ends up splitting runs even when the final formatting
# of all the pieces stays identical).
function Split-Replace($paragraph, $find, $replace) {
    $r = $paragraph.Range
    $text = $r.Text
    $idx = $text.IndexOf($find)
    $s = $r.Start + $idx
    $e = $s + $find.Length
    $sub = $d.Range($s, $e)
    $sub.Text = $replace
    $sub2 = $d.Range($s, $s + $replace.Length)
    $sub2.Font.Bold = $true
    $sub2.Font.Bold = $false
}

# --- Change 1 ---------------------------------------------------------
# "Poner un mensaje de error cuando el rut esta duplicado en la base de
#  datos" -> fix "esta" to "este" and split the sentence into 3 runs
#  around "este duplicado".
$p16 = $d.Paragraphs.Item(16)
Split-Replace $p16 "esta duplicado" "este duplicado"

# --- Change 2 ---------------------------------------------------------
# Insert a new bulleted requirement "Comentar bloques de codigo" right
# after the "Que sea responsive" bullet.
$p20 = $d.Paragraphs.Item(20)
$p20.Range.InsertParagraphAfter()
$p21 = $d.Paragraphs.Item(21)
$p21.Range.InsertAfter("Comentar bloques de codigo")

# --- Change 3 ---------------------------------------------------------
# "Tratar de "desglosar" ... "registrar un producto" este listo para su
#  uso..." -> fix "este" to "esté" and split the sentence into 3 runs
#  around "esté" (paragraph shifted from 29 to 30 after the insertion
#  above).
$p30 = $d.Paragraphs.Item(30)
Split-Replace $p30 "este" "esté"
